# "c1 modificada para incluir todas las paradas"
# Update the "parameters" fleet_size value and rework the
# "comp_quantity_inst1" sheet so it only keeps the origin/destination pair
# that covers every stop (collapsing the 3 rows down to 1).

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("parameters")
$wsComp   = $wb.Worksheets.Item("comp_quantity_inst1")

# --- parameters: fleet_size (B3) 2 -> 1 ---
$wsParams.Range("B3").Value = 1

# --- comp_quantity_inst1: keep a single row (old row 4 -> new row 2) ---
# overwrite row 2 with what used to be row 4 (F1 -> T1, qty 26, required)
$wsComp.Range("A2").Value = "F1"
$wsComp.Range("B2").Value = "T1"
$wsComp.Range("C2").Value = 26
$wsComp.Range("D2").Value = 1

# drop the now-redundant rows (old rows 3 and 4)
$wsComp.Rows.Item(3).Delete()
$wsComp.Rows.Item(3).Delete()

# reselect and hand focus back to the parameters sheet
$wsComp.Range("F20").Select()

$wsParams.Activate()
$wsParams.Range("B4").Select()
